$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.359.13"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.609.40"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Formula = "'607.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Formula = "'141.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("D7").Value = "3.608.49"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").Formula = "'7.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.30%  "
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("D13").Value = "4.230.77"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").Formula = "'28.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.66%  "
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "3.623.35"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.420.98"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Formula = "'0.117"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Formula = "'10.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Formula = "'14.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.41%  "
$ws.Range("D21").Formula = "'5.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Formula = "'399.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "3.755.85"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").Formula = "'75.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("D28").Formula = "'8.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").Formula = "'1.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +27.44%  "
$ws.Range("D30").Formula = "'8.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.99%  "
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "3.620.60"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Formula = "'24.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D37").Formula = "'5.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.17%  "
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").Formula = "'7.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").Formula = "'168.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Formula = "'0.0846"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").Formula = "'0.847"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Formula = "'26.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Formula = "'1.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.35%  "
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").Formula = "'43.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("D48").Formula = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Formula = "'7.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").Value = "2.475.38"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").Formula = "'0.920"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.34%  "
